$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$desc = "We are seeking a Senior RPA Developer to design, develop, and support automation solutions.`nCollaborate with teams to streamline business processes using RPA tools like UiPath or Automation Anywhere. Join Akkodis to grow your skills in a dynamic, tech-driven environment"

$ws.Range("A3").Value = "JD_002"
$ws.Range("B3").Value = "Senior RPA Developer"
$ws.Range("C3").Value = $desc
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 4
$ws.Range("F3").Value = "Remote"
$ws.Range("G3").Value = "Bengaluru, Karnataka, India"

$ws.Rows("3").AutoFit()
